$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.201571333333333
$ws.Range("H2").Value = 18.604714
$ws.Range("I2").Value = 0.05221490529364391
$ws.Range("J2").Value = 0.07406232529850043
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.768073999999999
$ws.Range("N2").Value = 14.304222
$ws.Range("O2").Value = 0.5213599352042468
$ws.Range("P2").Value = 0.5361472681932758
$ws.Range("Q2").Value = 29.56955103361199
$ws.Range("R2").Value = 266.125959302508
$ws.Range("S2").Value = 0.02722275964059007
$ws.Range("T2").Value = 0.03970831338483275

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.201571333333333
$ws.Range("H3").Value = 18.604714
$ws.Range("I3").Value = 0.05221490529364391
$ws.Range("J3").Value = 0.07406232529850043
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.082852
$ws.Range("N3").Value = 9.248556000000001
$ws.Range("O3").Value = 0.3370911439219029
$ws.Range("P3").Value = 0.3466520607784562
$ws.Range("Q3").Value = 19.11852658810933
$ws.Range("R3").Value = 172.066739292984
$ws.Range("S3").Value = 0.01760118215520825
$ws.Range("T3").Value = 0.02567385769076957

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.201571333333333
$ws.Range("H4").Value = 18.604714
$ws.Range("I4").Value = 0.05221490529364391
$ws.Range("J4").Value = 0.07406232529850043
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1813516666666667
$ws.Range("N4").Value = 0.544055
$ws.Range("O4").Value = 0.01982970339439268
$ws.Range("P4").Value = 0.02039213331538707
$ws.Range("Q4").Value = 1.124665297252222
$ws.Range("R4").Value = 10.12198767527
$ws.Range("S4").Value = 0.001035406084739263
$ws.Range("T4").Value = 0.001510288811134585

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.201571333333333
$ws.Range("H5").Value = 18.604714
$ws.Range("I5").Value = 0.05221490529364391
$ws.Range("J5").Value = 0.07406232529850043
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.3564626666666666
$ws.Range("N5").Value = 1.069388
$ws.Range("O5").Value = 0.0389770277885927
$ws.Range("P5").Value = 0.04008253331349799
$ws.Range("Q5").Value = 2.210628655003555
$ws.Range("R5").Value = 19.895657895032
$ws.Range("S5").Value = 0.002035181814609095
$ws.Range("T5").Value = 0.002968605621052269

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.201571333333333
$ws.Range("H6").Value = 18.604714
$ws.Range("I6").Value = 0.05221490529364391
$ws.Range("J6").Value = 0.07406232529850043
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 0.756715
$ws.Range("N6").Value = 1.51343
$ws.Range("O6").Value = 0.08274218969086504
$ws.Range("P6").Value = 0.05672600439938289
$ws.Range("Q6").Value = 4.692822051503333
$ws.Range("R6").Value = 28.15693230902
$ws.Range("S6").Value = 0.004320375598497237
$ws.Range("T6").Value = 0.004201259790711263

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 7.461641333333333
$ws.Range("H7").Value = 22.384924
$ws.Range("I7").Value = 0.06282422221945559
$ws.Range("J7").Value = 0.0891107233935555
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.768073999999999
$ws.Range("N7").Value = 14.304222
$ws.Range("O7").Value = 0.5213599352042468
$ws.Range("P7").Value = 0.5361472681932758
$ws.Range("Q7").Value = 35.57765803879199
$ws.Range("R7").Value = 320.1989223491279
$ws.Range("S7").Value = 0.03275403242559257
$ws.Range("T7").Value = 0.04777647091418142

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.461641333333333
$ws.Range("H8").Value = 22.384924
$ws.Range("I8").Value = 0.06282422221945559
$ws.Range("J8").Value = 0.0891107233935555
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.082852
$ws.Range("N8").Value = 9.248556000000001
$ws.Range("O8").Value = 0.3370911439219029
$ws.Range("P8").Value = 0.3466520607784562
$ws.Range("Q8").Value = 23.00313590774934
$ws.Range("R8").Value = 207.028223169744
$ws.Range("S8").Value = 0.02117748893396011
$ws.Range("T8").Value = 0.030890415901835

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.461641333333333
$ws.Range("H9").Value = 22.384924
$ws.Range("I9").Value = 0.06282422221945559
$ws.Range("J9").Value = 0.0891107233935555
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.1813516666666667
$ws.Range("N9").Value = 0.544055
$ws.Range("O9").Value = 0.01982970339439268
$ws.Range("P9").Value = 0.02039213331538707
$ws.Range("Q9").Value = 1.353181091868889
$ws.Range("R9").Value = 12.17862982682
$ws.Range("S9").Value = 0.001245785692595219
$ws.Range("T9").Value = 0.001817157751271965

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.461641333333333
$ws.Range("H10").Value = 22.384924
$ws.Range("I10").Value = 0.06282422221945559
$ws.Range("J10").Value = 0.0891107233935555
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.3564626666666666
$ws.Range("N10").Value = 1.069388
$ws.Range("O10").Value = 0.0389770277885927
$ws.Range("P10").Value = 0.04008253331349799
$ws.Range("Q10").Value = 2.659796567390222
$ws.Range("R10").Value = 23.938169106512
$ws.Range("S10").Value = 0.002448701455244444
$ws.Range("T10").Value = 0.003571783539012093

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 7.461641333333333
$ws.Range("H11").Value = 22.384924
$ws.Range("I11").Value = 0.06282422221945559
$ws.Range("J11").Value = 0.0891107233935555
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 0.756715
$ws.Range("N11").Value = 1.51343
$ws.Range("O11").Value = 0.08274218969086504
$ws.Range("P11").Value = 0.05672600439938289
$ws.Range("Q11").Value = 5.646335921553333
$ws.Range("R11").Value = 33.87801552932
$ws.Range("S11").Value = 0.005198213712063253
$ws.Range("T11").Value = 0.005054895287255021

$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 105.106922
$ws.Range("H12").Value = 210.213844
$ws.Range("I12").Value = 0.8849608724869005
$ws.Range("J12").Value = 0.836826951307944
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.768073999999999
$ws.Range("N12").Value = 14.304222
$ws.Range("O12").Value = 0.5213599352042468
$ws.Range("P12").Value = 0.5361472681932758
$ws.Range("Q12").Value = 501.1575820082279
$ws.Range("R12").Value = 3006.945492049368
$ws.Range("S12").Value = 0.4613831431380642
$ws.Range("T12").Value = 0.4486624838942617

$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 105.106922
$ws.Range("H13").Value = 210.213844
$ws.Range("I13").Value = 0.8849608724869005
$ws.Range("J13").Value = 0.836826951307944
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.082852
$ws.Range("N13").Value = 9.248556000000001
$ws.Range("O13").Value = 0.3370911439219029
$ws.Range("P13").Value = 0.3466520607784562
$ws.Range("Q13").Value = 324.029084701544
$ws.Range("R13").Value = 1944.174508209264
$ws.Range("S13").Value = 0.2983124728327345
$ws.Range("T13").Value = 0.2900877871858517

$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 105.106922
$ws.Range("H14").Value = 210.213844
$ws.Range("I14").Value = 0.8849608724869005
$ws.Range("J14").Value = 0.836826951307944
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.1813516666666667
$ws.Range("N14").Value = 0.544055
$ws.Range("O14").Value = 0.01982970339439268
$ws.Range("P14").Value = 0.02039213331538707
$ws.Range("Q14").Value = 19.06131548290333
$ws.Range("R14").Value = 114.36789289742
$ws.Range("S14").Value = 0.0175485116170582
$ws.Range("T14").Value = 0.01706468675298052

$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 105.106922
$ws.Range("H15").Value = 210.213844
$ws.Range("I15").Value = 0.8849608724869005
$ws.Range("J15").Value = 0.836826951307944
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.3564626666666666
$ws.Range("N15").Value = 1.069388
$ws.Range("O15").Value = 0.0389770277885927
$ws.Range("P15").Value = 0.04008253331349799
$ws.Range("Q15").Value = 37.46669370124533
$ws.Range("R15").Value = 224.800162207472
$ws.Range("S15").Value = 0.03449314451873917
$ws.Range("T15").Value = 0.03354214415343363

$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 105.106922
$ws.Range("H16").Value = 210.213844
$ws.Range("I16").Value = 0.8849608724869005
$ws.Range("J16").Value = 0.836826951307944
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 0.756715
$ws.Range("N16").Value = 1.51343
$ws.Range("O16").Value = 0.08274218969086504
$ws.Range("P16").Value = 0.05672600439938289
$ws.Range("Q16").Value = 79.53598448123
$ws.Range("R16").Value = 318.14393792492
$ws.Range("S16").Value = 0.07322360038030455
$ws.Range("T16").Value = 0.0474698493214166
